# Update the "Förändrad" (Changed) date column C for all data rows (C2:C495)
# from 45180 (2023-09-11) to 45181 (2023-09-12), leaving everything else untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C2:C495")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
